# Apply "hybrid bold + color" highlighting to quantitative impact metrics
# across several bullet paragraphs in the resume document.
#
# For each target paragraph we locate it via a stable (pre-edit) substring,
# then for every metric token inside that paragraph we use Range.Find to
# seek to that exact occurrence and flip Font.Bold / Font.Color on it.
# Word's Find.Execute naturally splits the enclosing run into a plain
# "before" run, the matched/now-formatted run, and a plain "after" run --
# exactly mirroring the <w:r> splits shown in the target diff.

$d = $word.ActiveDocument

# RGB(0x2C, 0x3E, 0x50) using Word's r + g*256 + b*65536 packing -> "2C3E50"
$metricColor = 5258796

function Set-MetricBold {
    param(
        [object]$Paragraph,
        [string]$Token
    )
    $rng = $Paragraph.Range
    $found = $rng.Find.Execute($Token, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: '$Token'"
        return
    }
    $rng.Font.Bold = $true
    $rng.Font.Color = $metricColor
}

function Find-ParagraphByText {
    param(
        [string]$Needle
    )
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$Needle*") {
            return $p
        }
    }
    return $null
}

function Find-ParagraphByExactText {
    # Matches a paragraph whose trimmed text equals $Needle exactly -- used
    # to disambiguate paragraphs whose text is a strict prefix/substring of
    # another paragraph's text (e.g. the short- vs. long-form "Achieved 87%"
    # bullets).
    param(
        [string]$Needle
    )
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Trim() -eq $Needle) {
            return $p
        }
    }
    return $null
}

# 1) "• Discovered systematic race coding errors ... from 23% to 64%"
$p1 = Find-ParagraphByText("Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning")
Set-MetricBold $p1 "23%"
Set-MetricBold $p1 "64%"

# 2) "• Achieved 87% ... industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"
$p2 = Find-ParagraphByText("reducing polling error margins")
Set-MetricBold $p2 "87%"
Set-MetricBold $p2 "71%"
Set-MetricBold $p2 "±4.2%"
Set-MetricBold $p2 "±2.1%"

# 3) "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
$p3 = Find-ParagraphByText("Wrote RFP and analyzed bids from")
Set-MetricBold $p3 "1,200"

# 4) "• Created comprehensive meta-analysis framework ... $400M ... now valued at $1B+"
$p4 = Find-ParagraphByText("Created comprehensive meta-analysis framework")
Set-MetricBold $p4 "$400M"
Set-MetricBold $p4 "$1B"

# 5) "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
$p5 = Find-ParagraphByText("Algorithm reduced mapping costs by")
Set-MetricBold $p5 "73.5%"
Set-MetricBold $p5 "$4.7M"

# 6) "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (no margins clause)
$p6 = Find-ParagraphByExactText("• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%")
Set-MetricBold $p6 "87%"
Set-MetricBold $p6 "71%"

Write-Host "Done"
